$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Scenathon2023 -> SOFA 2024" data refresh added two alpha3/country
# mapping rows that already existed elsewhere in the table (GRC/Greece and
# NPL/Nepal each now appear twice), shifting every row below them down by one.

# Insert GRC / Greece as the new row 11 (existing data shifts down).
$ws.Rows(11).Insert()
$ws.Range("A11").Value = "GRC"
$ws.Range("B11").Value = "Greece"

# Insert NPL / Nepal as the new row 24 (after NOR / Norway, existing data
# shifts down again).
$ws.Rows(24).Insert()
$ws.Range("A24").Value = "NPL"
$ws.Range("B24").Value = "Nepal"

# Widen column A so the longer region codes/names are readable, and leave the
# cursor where the author's session ended up (A22).
$ws.Columns("A").ColumnWidth = 20.8
$ws.Range("A22").Select() | Out-Null
